$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental = true
# (use formula + copy/paste-values so "true" lands as text, not a Boolean)
$ws.Range("B7").Formula = '="true"'
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null

# Date updated
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"

# Compositional = false
$ws.Range("B18").Formula = '="false"'
$ws.Range("B18").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0
